$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 6 (B6:D6) to the new consolidated text, matching B6's style.
$newText = "60% of the respondents agree to use the system"
$ws.Range("B6").Value = $newText
$ws.Range("C6").Value = $newText
$ws.Range("D6").Value = $newText

# Make C6/D6 share B6's style, and grow the row height to fit wrapped text.
$ws.Range("C6:D6").Style = $ws.Range("B6").Style
$ws.Rows.Item(6).RowHeight = 45

# Update the view: scroll so row 5 is the top-left visible row, and select B6:D6.
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("B6:D6").Select()
